$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '27.136.99'
$ws.Range("E2").Value = '  +0.44%  '
$ws.Range("D3").Value = '1.678.13'
$ws.Range("E3").Value = '  -0.06%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '214.24'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -0.73%  '
$ws.Range("E6").Value = '  -0.08%  '
$ws.Range("E7").Value = '  +0.08%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '22.79'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  +6.55%  '
$ws.Range("E9").Value = '  +2.18%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.0622'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  -0.47%  '
$ws.Range("E11").Value = '  +0.17%  '
$ws.Range("D12").Value = '1.914.64'
$ws.Range("E12").Value = '  -0.05%  '
$ws.Range("D13").Value = '1.680.25'
$ws.Range("E13").Value = '  +0.19%  '
$ws.Range("E14").Value = '  +2.21%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.552'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  +3.31%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '66.54'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  +0.06%  '
$ws.Range("D17").Value = '27.101.36'
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '234.96'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  -0.47%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '7.89'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  -3.28%  '
$ws.Range("E20").Value = '  +0.36%  '
$ws.Range("E21").Value = '  +0.05%  '
$ws.Range("E22").Value = '  +1.49%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '9.52'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  +2.68%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '2.09'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  -1.33%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '148.32'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  +1.25%  '
$ws.Range("E26").Value = '  +2.20%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '16.34'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  -0.61%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '0.112'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  -0.25%  '
$ws.Range("E29").Value = '  +0.27%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '0.0501'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  +0.52%  '
$ws.Range("E31").Value = '  -0.47%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '3.36'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  -0.33%  '
$ws.Range("D33").Value = '1.539.86'
$ws.Range("E33").Value = '  -0.11%  '
$ws.Range("E34").Value = '  +1.26%  '
$ws.Range("E35").Value = '  -3.58%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.605'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  +2.67%  '
$ws.Range("E37").Value = '  +2.05%  '
$ws.Range("E38").Value = '  -0.15%  '
$ws.Range("E39").Value = '  -1.18%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '1.06'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  +2.22%  '
$ws.Range("E41").Value = '  +3.26%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '69.39'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  +2.12%  '
$ws.Range("E43").Value = '  +0.09%  '
$ws.Range("E44").Value = '  -0.29%  '
$ws.Range("D45").Value = '1.822.05'
$ws.Range("E45").Value = '  +0.12%  '
$ws.Range("E46").Value = '  -0.24%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '89.89'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  -0.53%  '
$ws.Range("E48").Value = '  +5.93%  '
$ws.Range("E49").Value = '  +2.89%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '8.18'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  +2.45%  '
$ws.Range("E51").Value = '  -0.35%  '
